$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 682, pushing existing rows 682:734 down to 683:735
$ws.Rows(682).Insert()

# Populate the newly inserted row 682 with the new data point
$ws.Cells.Item(682, 1).Value = 3
$ws.Cells.Item(682, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(682, 3).Value = "Coquimbo"
$ws.Cells.Item(682, 4).Value = 45223
$ws.Cells.Item(682, 5).Value = 5
$ws.Cells.Item(682, 6).Value = 100112017
$ws.Cells.Item(682, 7).Value = "Apio"
$ws.Cells.Item(682, 8).Value = "Americana (o)"
$ws.Cells.Item(682, 9).Value = "Primera"
$ws.Cells.Item(682, 10).Value = 250
$ws.Cells.Item(682, 11).Value = 8000
$ws.Cells.Item(682, 12).Value = 9000
$ws.Cells.Item(682, 13).Value = 8400
$ws.Cells.Item(682, 14).Value = "`$/docena de matas"
$ws.Cells.Item(682, 15).Value = "Pan de Azúcar"
$ws.Cells.Item(682, 16).Value = 1400
$ws.Cells.Item(682, 17).Value = 6
$ws.Cells.Item(682, 18).Value = "Hortaliza"
